$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet2 ("Semilla 9") updates ---

# New cells in row 4 (host ssh credentials)
$ws2.Range("E4").Value = "host ssh"
$ws2.Range("F4").Value = "usuario ssh"
$ws2.Range("G4").Value = "contraseña ssh"

# New cells in row 5
$ws2.Range("E5").Value = "10.69.60.119"
$ws2.Range("F5").Value = "consulta_log"
$ws2.Range("G5").Value = "consulta_log"

# New cells in row 8 (header row for the new block)
$ws2.Range("E8").Value = "msisdn"
$ws2.Range("F8").Value = "msi"
$ws2.Range("G8").Value = "nip"
$ws2.Range("H8").Value = "client"
$ws2.Range("I8").Value = "portId"

# Row 9: update existing C/D and add new E..I
$ws2.Range("C9").Value = "3043208091"
$ws2.Range("D9").Value = "732111324707274"
$ws2.Range("E9").Value = "3045981684"
$ws2.Range("F9").Value = "732111193278813"
$ws2.Range("G9").Value = "81684"
$ws2.Range("H9").Value = "1061520830"
$ws2.Range("I9").Value = "00002201108240181684"

# Row 10: update existing C/D and add new E..I
$ws2.Range("C10").Value = "3043209773"
$ws2.Range("D10").Value = "732111324707275"
$ws2.Range("E10").Value = "3045984642"
$ws2.Range("F10").Value = "732111193278730"
$ws2.Range("G10").Value = "81670"
$ws2.Range("H10").Value = "111295346"
$ws2.Range("I10").Value = "00002201108240181670"

# Row 11: update existing C/D
$ws2.Range("C11").Value = "3043208091"
$ws2.Range("D11").Value = "732111324707274"

# Row 12: update existing C/D
$ws2.Range("C12").Value = "3043209819"
$ws2.Range("D12").Value = "732111324707276"

# Row 13: update existing C/D
$ws2.Range("C13").Value = "3043209868"
$ws2.Range("D13").Value = "732111324707278"

# Row 14: brand new row, mirroring row 13's A/B and new C/D values.
# Force text number format first so numeric-looking values stay text (t="s"),
# matching the rest of the sheet.
$ws2.Range("A14").NumberFormat = "@"
$ws2.Range("A14").Value = "10960370"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "1050388676"
$ws2.Range("C14").NumberFormat = $ws3.Range("C14").NumberFormat
$ws2.Range("C14").HorizontalAlignment = $ws3.Range("C14").HorizontalAlignment
$ws2.Range("C14").Value = "3043209863"
$ws2.Range("D14").NumberFormat = "@"
$ws2.Range("D14").Value = "732111324707277"

# Update selection for sheet2 (no tab activation should remain on it)
$ws2.Range("D18").Select()

# --- Sheet3 ("Semilla 8") updates ---
# Swap C13 / C14 values
$ws3.Range("C13").Value = "3043209868"
$ws3.Range("C14").Value = "3043209863"

# Restore sheet3 as the active/selected tab with new selection
$ws3.Range("B15").Select()
